$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 146, shifting rows 146:260 down to 147:261.
$ws.Rows.Item(146).Insert()

# New row 146 duplicates the (now-shifted) original row 146 content (now at row 147),
# except for the Fecha (D) and Volumen (J) values which change.
$srcRow = 147
$dstRow = 146

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item($dstRow, $col).Value = $ws.Cells.Item($srcRow, $col).Value2
}

# Copy the date cell's number format/style too
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat

# Now set the two changed values
$ws.Cells.Item($dstRow, 4).Value = 44729
$ws.Cells.Item($dstRow, 10).Value = 180
